$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value = 44595
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100103
$ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value = 100103002
$ws.Cells.Item($row, 10).Value = "Ciruela"
$ws.Cells.Item($row, 11).Value = "Black Amber"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 60
$ws.Cells.Item($row, 14).Value = 9000
$ws.Cells.Item($row, 15).Value = 9500
$ws.Cells.Item($row, 16).Value = 9250
$ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row, 19).Value = 514
$ws.Cells.Item($row, 20).Value = 18
